# feat(database): Integrate MongoDB Atlas for persistent data storage
#
# Adds a "timestamp" column (H) to the reorder-optimization
# recommendations sheet, recording the moment each recommendation
# document was persisted (mirrors the timestamp MongoDB Atlas stamps
# on inserted records).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Serial date value for the persistence timestamp applied to every row.
$timestamp = 45764.31028708907

# --- Header -----------------------------------------------------------
$ws.Range("H1").Value = "timestamp"
# Match the look of the existing header cells (bold font, borders,
# centered alignment) by copying their formatting onto the new header.
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows ----------------------------------------------------------
$lastRow = $ws.Cells(1, 1).End(-4121).Row

# First data row establishes the custom datetime number format (applied
# once in lower case, then normalized to upper case) used for the whole
# timestamp column.
$firstCell = $ws.Cells.Item(2, 8)
$firstCell.Value = $timestamp
$firstCell.NumberFormat = "yyyy-mm-dd h:mm:ss"
$firstCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($r = 3; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Value = $timestamp
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
